# Generate Report for Handback
#
# A new handback event is recorded for e6ec6f30-8f85-493d-bf8f-2306185dcdff.md
# on both the "zh-cn" and "de-de" worksheets (row 7 of each status table).
# The handback turned out to not be based on the latest source, so:
#   - Latest Target File (I7) gets a hyperlink to the .md file (pointing at
#     the latest commit of the source).
#   - Latest Handback File (J7) is set to the locale-specific xlf file that
#     was just handed back (same value as the existing Latest Handoff File
#     in column G).
#   - Latest Handback DateTime (K7) records when this handback happened.
#   - Error Detail (P7) records that the handback was not based on the
#     latest version of the source file.

$wb = $excel.ActiveWorkbook

$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/572094d48d71b4964ddcb49cb22935899c4cf1d6/e2e/e6ec6f30-8f85-493d-bf8f-2306185dcdff.md"
$notLatestMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7fe0ce6d18e868f315e6ab3df8ba4a215ced088/e2e/e6ec6f30-8f85-493d-bf8f-2306185dcdff.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/572094d48d71b4964ddcb49cb22935899c4cf1d6/e2e/e6ec6f30-8f85-493d-bf8f-2306185dcdff.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $latestMdUrl, [Type]::Missing, [Type]::Missing, "e6ec6f30-8f85-493d-bf8f-2306185dcdff.md")
$wsZh.Range("J7").Value = $wsZh.Range("G7").Value()
$wsZh.Range("K7").Value = "2016-08-22 06:53:55"
$wsZh.Range("P7").Value = $notLatestMessage

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $latestMdUrl, [Type]::Missing, [Type]::Missing, "e6ec6f30-8f85-493d-bf8f-2306185dcdff.md")
$wsDe.Range("J7").Value = $wsDe.Range("G7").Value()
$wsDe.Range("K7").Value = "2016-08-22 06:54:05"
$wsDe.Range("P7").Value = $notLatestMessage
